$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 129, shifting existing rows 129:228 down to 130:228
$ws.Rows.Item(129).Insert()

# Populate the newly inserted row 129 with the new data record
$ws.Cells.Item(129, 1).Value = 3
$ws.Cells.Item(129, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(129, 3).Value = "Coquimbo"
$ws.Cells.Item(129, 4).Value = 44942
$ws.Cells.Item(129, 5).Value = 5
$ws.Cells.Item(129, 6).Value = 100112030
$ws.Cells.Item(129, 7).Value = "Poroto granado"
$ws.Cells.Item(129, 8).Value = "Sin especificar"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 73
$ws.Cells.Item(129, 11).Value = 44000
$ws.Cells.Item(129, 12).Value = 45000
$ws.Cells.Item(129, 13).Value = 44479
$ws.Cells.Item(129, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(129, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(129, 16).Value = 1779
$ws.Cells.Item(129, 17).Value = 25
$ws.Cells.Item(129, 18).Value = "Hortaliza"
